# The edit described by the diff is a pure re-ordering of the data rows
# 19-22, 24-27, 29-30 and 32-35 (rows 23, 28 and 31 keep their content):
# the cell values themselves don't change, entire rows just swap places
# with each other, grouped into the following permutation cycles:
#   19 -> 20 -> 24 -> 19
#   21 -> 25 -> 22 -> 26 -> 21
#   27 -> 29 -> 30 -> 27
#   32 -> 34 -> 33 -> 35 -> 32
# For every cycle we stash each row's full contents (columns A:AY) in a
# scratch row far below the data, wipe the source rows (so no stray
# value survives in a column that the destination row didn't already
# use) and then copy the stashed rows back into their new home.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cycle: [19, 20, 24]
$ws.Range("A19:AY19").Copy($ws.Range("A100:AY100"))
$ws.Range("A20:AY20").Copy($ws.Range("A101:AY101"))
$ws.Range("A24:AY24").Copy($ws.Range("A102:AY102"))
$ws.Range("A19:AY19").Clear()
$ws.Range("A20:AY20").Clear()
$ws.Range("A24:AY24").Clear()
$ws.Range("A100:AY100").Copy($ws.Range("A20:AY20"))
$ws.Range("A101:AY101").Copy($ws.Range("A24:AY24"))
$ws.Range("A102:AY102").Copy($ws.Range("A19:AY19"))
$ws.Range("A100:AY102").Clear()

# Cycle: [21, 25, 22, 26]
$ws.Range("A21:AY21").Copy($ws.Range("A100:AY100"))
$ws.Range("A25:AY25").Copy($ws.Range("A101:AY101"))
$ws.Range("A22:AY22").Copy($ws.Range("A102:AY102"))
$ws.Range("A26:AY26").Copy($ws.Range("A103:AY103"))
$ws.Range("A21:AY21").Clear()
$ws.Range("A25:AY25").Clear()
$ws.Range("A22:AY22").Clear()
$ws.Range("A26:AY26").Clear()
$ws.Range("A100:AY100").Copy($ws.Range("A25:AY25"))
$ws.Range("A101:AY101").Copy($ws.Range("A22:AY22"))
$ws.Range("A102:AY102").Copy($ws.Range("A26:AY26"))
$ws.Range("A103:AY103").Copy($ws.Range("A21:AY21"))
$ws.Range("A100:AY103").Clear()

# Cycle: [27, 29, 30]
$ws.Range("A27:AY27").Copy($ws.Range("A100:AY100"))
$ws.Range("A29:AY29").Copy($ws.Range("A101:AY101"))
$ws.Range("A30:AY30").Copy($ws.Range("A102:AY102"))
$ws.Range("A27:AY27").Clear()
$ws.Range("A29:AY29").Clear()
$ws.Range("A30:AY30").Clear()
$ws.Range("A100:AY100").Copy($ws.Range("A29:AY29"))
$ws.Range("A101:AY101").Copy($ws.Range("A30:AY30"))
$ws.Range("A102:AY102").Copy($ws.Range("A27:AY27"))
$ws.Range("A100:AY102").Clear()

# Cycle: [32, 34, 33, 35]
$ws.Range("A32:AY32").Copy($ws.Range("A100:AY100"))
$ws.Range("A34:AY34").Copy($ws.Range("A101:AY101"))
$ws.Range("A33:AY33").Copy($ws.Range("A102:AY102"))
$ws.Range("A35:AY35").Copy($ws.Range("A103:AY103"))
$ws.Range("A32:AY32").Clear()
$ws.Range("A34:AY34").Clear()
$ws.Range("A33:AY33").Clear()
$ws.Range("A35:AY35").Clear()
$ws.Range("A100:AY100").Copy($ws.Range("A34:AY34"))
$ws.Range("A101:AY101").Copy($ws.Range("A33:AY33"))
$ws.Range("A102:AY102").Copy($ws.Range("A35:AY35"))
$ws.Range("A103:AY103").Copy($ws.Range("A32:AY32"))
$ws.Range("A100:AY103").Clear()
